$wb = $excel.ActiveWorkbook
$wsCompleteness = $wb.Worksheets.Item("Completeness")
$wsInstructions = $wb.Worksheets.Item("Instructions")

$wsInstructions.Activate()
$wsInstructions.Range("C1").Value = "Template updated 12/8/22."
$wsInstructions.Range("C1").Font.Color = 255
$wsInstructions.Range("C3").Select()

$wsCompleteness.Activate()
